$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1358.4082
$ws.Range("I15").Value = 1358.4082
$ws.Range("K15").Value = 4075.2246
$ws.Range("M15").Value = -3906.2246

$ws.Range("H18").Value = 83338740
$ws.Range("I18").Value = 83338740
$ws.Range("K18").Value = 83338740
$ws.Range("M18").Value = -83338456

$ws.Range("H28").Value = 1961.3182
$ws.Range("I28").Value = 1646
$ws.Range("J28").Value = 2637
$ws.Range("K28").Value = 1646
$ws.Range("L28").Value = 2637
$ws.Range("M28").Value = -1161
$ws.Range("N28").Value = -3607

$ws.Range("H40").Value = 25885
$ws.Range("I40").Value = 21499.8
$ws.Range("J40").Value = 28321.223
$ws.Range("K40").Value = 21499.8
$ws.Range("L40").Value = 28321.223
$ws.Range("M40").Value = -21324.8
$ws.Range("N40").Value = -28671.223

$ws.Range("H62").Value = 3933.8333
$ws.Range("I62").Value = 2701.6667
$ws.Range("K62").Value = 2701.6667
$ws.Range("M62").Value = -2077.6667

$ws.Range("H65").Value = 3933.8333
$ws.Range("I65").Value = 2701.6667
$ws.Range("K65").Value = 13508.3335
$ws.Range("M65").Value = -10388.3335

$ws.Range("H132").Value = 7910.705
$ws.Range("I132").Value = 3969.1968
$ws.Range("K132").Value = 11907.5904
$ws.Range("M132").Value = -9377.590400000001

$ws.Range("H137").Value = 14495173
$ws.Range("I137").Value = 2156.8462
$ws.Range("K137").Value = 6470.5386
$ws.Range("M137").Value = -3920.5386

$ws.Range("H141").Value = 5592.3076
$ws.Range("I141").Value = 5895.8335
$ws.Range("J141").Value = 1950
$ws.Range("K141").Value = 17687.5005
$ws.Range("L141").Value = 5850
$ws.Range("M141").Value = -12507.5005
$ws.Range("N141").Value = -16210

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H32").Value = 14791.604
$ws.Range("I32").Value = 14703.63
$ws.Range("J32").Value = 15319.444
$ws.Range("K32").Value = 14703.63
$ws.Range("L32").Value = 15319.444
$ws.Range("M32").Value = -14416.63
$ws.Range("N32").Value = -15893.444

$ws.Range("H61").Value = 3698.6206
$ws.Range("I61").Value = 3260.875
$ws.Range("K61").Value = 3260.875
$ws.Range("M61").Value = -3048.875

$ws.Range("H74").Value = 17859050
$ws.Range("I74").Value = 25001294
$ws.Range("K74").Value = 25001294
$ws.Range("M74").Value = -25000420

$ws.Range("H77").Value = 17859050
$ws.Range("I77").Value = 25001294
$ws.Range("K77").Value = 125006470
$ws.Range("M77").Value = -125002102

$ws.Range("H97").Value = 1483.4166
$ws.Range("I97").Value = 688.875
$ws.Range("K97").Value = 688.875
$ws.Range("M97").Value = -192.875

$ws.Range("H132").Value = 22507.969
$ws.Range("I132").Value = 26676.652
$ws.Range("K132").Value = 80029.95599999999
$ws.Range("M132").Value = -77499.95599999999

$ws.Range("H136").Value = 3698.6206
$ws.Range("I136").Value = 3260.875
$ws.Range("K136").Value = 9782.625
$ws.Range("M136").Value = -7232.625

$ws.Range("H138").Value = 71904.836
$ws.Range("J138").Value = 71904.836
$ws.Range("L138").Value = 71904.836
$ws.Range("N138").Value = -82184.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1486.0385
$ws.Range("J20").Value = 1806.8572
$ws.Range("L20").Value = 1806.8572
$ws.Range("N20").Value = -2300.8572

$ws.Range("H86").Value = 2301.375
$ws.Range("I86").Value = 2038.7273
$ws.Range("J86").Value = 2879.2
$ws.Range("K86").Value = 2038.7273
$ws.Range("L86").Value = 2879.2
$ws.Range("M86").Value = -915.7273
$ws.Range("N86").Value = -5125.2

$ws.Range("H89").Value = 2301.375
$ws.Range("I89").Value = 2038.7273
$ws.Range("J89").Value = 2879.2
$ws.Range("K89").Value = 10193.6365
$ws.Range("L89").Value = 14396
$ws.Range("M89").Value = -4577.636500000001
$ws.Range("N89").Value = -25628

$ws.Range("H94").Value = 1142893.9
$ws.Range("I94").Value = 1957948.2
$ws.Range("J94").Value = 1817.6
$ws.Range("K94").Value = 1957948.2
$ws.Range("L94").Value = 1817.6
$ws.Range("M94").Value = -1957497.2
$ws.Range("N94").Value = -2719.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 39636.547
$ws.Range("I62").Value = 18000.285
$ws.Range("K62").Value = 18000.285
$ws.Range("M62").Value = -17376.285

$ws.Range("H65").Value = 39636.547
$ws.Range("I65").Value = 18000.285
$ws.Range("K65").Value = 90001.425
$ws.Range("M65").Value = -86881.425

$ws.Range("H105").Value = 2067817.5
$ws.Range("I105").Value = 2842332.5
$ws.Range("K105").Value = 2842332.5
$ws.Range("M105").Value = -2840585.5

$ws.Range("H132").Value = 34201412
$ws.Range("I132").Value = 41676040
$ws.Range("J132").Value = 31699.285
$ws.Range("K132").Value = 125028120
$ws.Range("L132").Value = 95097.855
$ws.Range("M132").Value = -125025590
$ws.Range("N132").Value = -100157.855

$ws.Range("H134").Value = 2199.25
$ws.Range("I134").Value = 2253.0386
$ws.Range("K134").Value = 6759.1158
$ws.Range("M134").Value = -4224.1158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 634
$ws.Range("J12").Value = 1290.125
$ws.Range("L12").Value = 3870.375
$ws.Range("N12").Value = -4216.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10139584
$ws.Range("J11").Value = 13334444
$ws.Range("L11").Value = 13334444
$ws.Range("N11").Value = -13334722

$ws.Range("H14").Value = 8865000
$ws.Range("I14").Value = 10631999
$ws.Range("J14").Value = 30000
$ws.Range("K14").Value = 10631999
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = -10631831
$ws.Range("N14").Value = -30336

$ws.Range("H70").Value = 3502432.5
$ws.Range("I70").Value = 6497518
$ws.Range("K70").Value = 6497518
$ws.Range("M70").Value = -6497248

$ws.Range("H73").Value = 3502432.5
$ws.Range("I73").Value = 6497518
$ws.Range("K73").Value = 6497518
$ws.Range("M73").Value = -6496582

$ws.Range("H96").Value = 17553.334
$ws.Range("J96").Value = 17553.334
$ws.Range("L96").Value = 17553.334
$ws.Range("N96").Value = -23045.334

$ws.Range("H97").Value = 862.6539
$ws.Range("I97").Value = 805.6842
$ws.Range("J97").Value = 1017.2857
$ws.Range("K97").Value = 805.6842
$ws.Range("L97").Value = 1017.2857
$ws.Range("M97").Value = -309.6842
$ws.Range("N97").Value = -2009.2857

$ws.Range("H138").Value = 75000
$ws.Range("I138").Value = 75000
$ws.Range("K138").Value = 75000
$ws.Range("M138").Value = -69860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 47619436
$ws.Range("I55").Value = 83333510
$ws.Range("K55").Value = 83333510
$ws.Range("M55").Value = -83333337

$ws.Range("H132").Value = 5872.875
$ws.Range("I132").Value = 5136.1
$ws.Range("K132").Value = 15408.3
$ws.Range("M132").Value = -12878.3

$ws.Range("H136").Value = 4294.5
$ws.Range("I136").Value = 2257.7
$ws.Range("K136").Value = 6773.099999999999
$ws.Range("M136").Value = -4223.099999999999

$ws.Range("H139").Value = 100714
$ws.Range("J139").Value = 100714
$ws.Range("L139").Value = 100714
$ws.Range("N139").Value = -110994

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8763.125
$ws.Range("J74").Value = 8763.125
$ws.Range("L74").Value = 8763.125
$ws.Range("N74").Value = -10635.125

$ws.Range("H77").Value = 8763.125
$ws.Range("J77").Value = 8763.125
$ws.Range("L77").Value = 26289.375
$ws.Range("N77").Value = -35649.375

$ws.Range("H132").Value = 45459384
$ws.Range("I132").Value = 2003
$ws.Range("J132").Value = 62505900
$ws.Range("K132").Value = 6009
$ws.Range("L132").Value = 187517700
$ws.Range("M132").Value = -3479
$ws.Range("N132").Value = -187522760

$ws.Range("H136").Value = 7176.19
$ws.Range("I136").Value = 3069.4348
$ws.Range("J136").Value = 10674.537
$ws.Range("K136").Value = 9208.304400000001
$ws.Range("L136").Value = 32023.611
$ws.Range("M136").Value = -6658.304400000001
$ws.Range("N136").Value = -37123.611
